$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workflow_repository_tools")

# Insert a new row right below the header (pushing ControlFreeC and
# everything after it down by one row). Copying an existing data row
# first (instead of a bare Rows.Insert) means the new row inherits the
# same "text" cell typing the rest of the table uses for every column,
# so plain numeric-looking strings like "1.3" don't silently become
# numbers when we overwrite them below.
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(2).Insert()
$excel.CutCopyMode = $false

$ws.Range("A2").Value = "NGSCheckMate"
$ws.Range("C2").Value = "https://github.com/d3b-center/OpenPBTA-workflows/blob/master/cwl/bcf_call.cwl"

# "1.3" looks numeric, so a plain .Value assignment would store it as a
# number instead of text (unlike the rest of this "Version" column).
# Compute it as a formula result (a text string) in a scratch cell, then
# paste only the resulting value into B2 so it keeps its text type -
# this also lets it reuse the existing shared "1.3" string instead of
# leaving behind a stray number-formatted style.
$ws.Range("Z1").Formula = '="1.3"'
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$excel.CutCopyMode = $false
